$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "41.839.53"
$ws.Range("D3").Value = "2.230.54"
$ws.Range("D5").Value = "'231.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'60.61"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Value = "'0.406"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'58.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Value = "'0.104"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "2.560.46"
$ws.Range("D14").Value = "'15.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'22.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("D18").Value = "2.252.20"
$ws.Range("D19").Value = "41.733.81"
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("D21").Value = "'72.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'6.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'248.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'2.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Value = "'9.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'169.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'0.142"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'19.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("D35").Value = "'4.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'0.0652"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'6.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'3.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'2.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").Value = "'0.0240"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'8.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Value = "'98.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.0959"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "1.472.07"
$ws.Range("D49").Value = "'16.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'2.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'2.79"
$ws.Range("D51").Style = "Normal"

# --- Column E (Volume/1h change) updates ---
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("E6").Value = "  -1.49%  "
$ws.Range("E7").Value = "  -6.59%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("E16").Value = "  -2.68%  "
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  -2.28%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("E28").Value = "  -2.21%  "
$ws.Range("E29").Value = "  -2.59%  "
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("E31").Value = "  -2.48%  "
$ws.Range("E32").Value = "  -8.14%  "
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("E34").Value = "  +3.97%  "
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("E36").Value = "  +3.11%  "
$ws.Range("E37").Value = "  -8.68%  "
$ws.Range("E38").Value = "  -5.01%  "
$ws.Range("E39").Value = "  -3.50%  "
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("E41").Value = "  +7.26%  "
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("E45").Value = "  -3.19%  "
$ws.Range("E46").Value = "  -8.37%  "
$ws.Range("E47").Value = "  +1.61%  "
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("E49").Value = "  -5.08%  "
$ws.Range("E50").Value = "  +8.99%  "
$ws.Range("E51").Value = "  -2.43%  "

# --- Rows 50/51: coin rank swap (NEARProtocol moves up, HuobiToken moves down) ---
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
